$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrap the existing long description (currently in G2) in <p> tags to make it HTML.
$ws.Range("G2").Value = "<p>A vesperák az egyházi év vasárnapjainak megfelelően változó anyagát a Gyülekezeti liturgikus könyv 2007-es megjelenéséig az ünnepkörönként változó „küllap” és az adott vasárnap jellegét tükröző „bellap” kombinálásával adhattuk a gyülekezeti tagok kezébe. A kezdeti időszakban szöveges magyarázat is segítette a tájékozódást.</p>"

# Add a new HTML description for row 3's "Leírás" (Description) column.
$ws.Range("G3").Value = "<p>Első bekezdés</p>`n<p>Második bekezdés</p>"

# Wrap text on the new cell and make the row taller to fit (2 lines @ 15pt).
$ws.Range("G3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 30

# Move the selection, matching the author's final cursor position.
$ws.Range("G13").Select()
